# Apply updated "想去人数" (F) / "最低票价" (G) values to the
# "展览" and "全部类型" worksheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

# Cell updates: row -> @{ column letter = new value }
$updates = @{
    2  = @{ G = 65 }
    3  = @{ F = 358; G = 50 }
    4  = @{ F = 1867; G = 60 }
    8  = @{ F = 743 }
    11 = @{ F = 4452 }
    13 = @{ F = 340 }
    14 = @{ F = 1237 }
    17 = @{ F = 810 }
    18 = @{ F = 27 }
    19 = @{ F = 438 }
    21 = @{ F = 215 }
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$row").Value = $cols[$col]
        }
    }
}
